$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.181780576705933
$ws.Range("B1").Value = 2.423226118087769
$ws.Range("D1").Value = 2.325670719146729
$ws.Range("E1").Value = 1.190549850463867
